# Re-knit: update the standard-error columns (n_cbh_mean_se, n_cbh_median_se)
# with freshly recomputed values for every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.172253206300272
$ws.Range("H2").Value = 0.476721368582596

$ws.Range("F3").Value = 0.199665905482454
$ws.Range("H3").Value = 0.469101328132906

$ws.Range("F4").Value = 0.248757258059664
$ws.Range("H4").Value = 0.534978693222639

$ws.Range("F5").Value = 0.267529967255851
$ws.Range("H5").Value = 0.498766245617892

$ws.Range("F6").Value = 0.309351349624615
$ws.Range("H6").Value = 0.49858257549235

$ws.Range("F7").Value = 0.28790591185566
$ws.Range("H7").Value = 0.399894881082411

$ws.Range("F8").Value = 0.287700673587536
$ws.Range("H8").Value = 0.4867922726551

$ws.Range("F9").Value = 0.280875507150911
$ws.Range("H9").Value = 0.473934979047089

$ws.Range("F10").Value = 0.475375923256966
$ws.Range("H10").Value = 0.544632329534682

$ws.Range("F11").Value = 0.182751306265823
$ws.Range("H11").Value = 0.0611022379580608

$ws.Range("F12").Value = 0.186327429501234
$ws.Range("H12").Value = 0.476113888092371

$ws.Range("F13").Value = 0.292241694552718
$ws.Range("H13").Value = 0.447991848917697

$ws.Range("F14").Value = 0.41872055168621
$ws.Range("H14").Value = 0.438645825422014

$ws.Range("F15").Value = 0.449126226176202
$ws.Range("H15").Value = 0.724030538182072

$ws.Range("F16").Value = 1.11476960150764
$ws.Range("H16").Value = 1.34602629498342

$ws.Range("F17").Value = 2.24031699427435
$ws.Range("H17").Value = 3.08725803297395

$ws.Range("F18").Value = 1.70037246091027
$ws.Range("H18").Value = 1.4235949643851

$ws.Range("F19").Value = 1.60199758409185
$ws.Range("H19").Value = 2.03003079410757

$ws.Range("F20").Value = 1.47293442595309
$ws.Range("H20").Value = 1.76977855002148
